# Commit: "Push Excel and Change email details"
#
# 1. addListItem!A2 ("Surats" -> "SuratT") - shared string edit, ripples into
#    the formula-driven C2 ("=A2").
# 2. createUser!A2 (65 -> 66) - ripples into the CONCAT-driven B2 (user name)
#    and F2 (email address), which is where the "Change email details" part
#    of the commit message comes from.
# 3. The active/selected sheet moves from "addListItem" to "createUser", with
#    the selection on createUser left at A2.

$wb = $excel.ActiveWorkbook

$wsAddListItem = $wb.Worksheets.Item("addListItem")
$wsAddListItem.Range("A2").Value = "SuratT"

$wsCreateUser = $wb.Worksheets.Item("createUser")
$wsCreateUser.Range("A2").Value = 66

# Make "createUser" the active sheet/tab, with A2 selected, matching the
# updated workbookView.activeTab / sheetView.tabSelected + selection in the
# target workbook.
$wsCreateUser.Activate()
$wsCreateUser.Range("A2").Select()
